# Daily attendance processing - 2025-10-09 05:20:21
# Reorders the "Recorded By" (column G) comma-separated author lists so that
# System / backup@backdoor.com / admin@admin.com entries are re-sequenced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, System, system"
    3   = "System, dnasr281@gmail.com"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    6   = "System, dnasr281@gmail.com"
    11  = "System, dnasr281@gmail.com"
    12  = "System, dnasr281@gmail.com"
    13  = "System, dnasr281@gmail.com"
    14  = "System, dnasr281@gmail.com"
    15  = "System, dnasr281@gmail.com"
    29  = "backup@backdoor.com, System, system"
    30  = "System, dnasr281@gmail.com"
    32  = "backup@backdoor.com, System"
    33  = "System, dnasr281@gmail.com"
    38  = "System, dnasr281@gmail.com"
    39  = "System, dnasr281@gmail.com"
    40  = "System, dnasr281@gmail.com"
    41  = "System, dnasr281@gmail.com"
    42  = "System, dnasr281@gmail.com"
    56  = "backup@backdoor.com, System, system"
    57  = "System, dnasr281@gmail.com"
    58  = "backup@backdoor.com, System"
    59  = "backup@backdoor.com, System"
    60  = "System, dnasr281@gmail.com"
    65  = "System, dnasr281@gmail.com"
    66  = "System, dnasr281@gmail.com"
    67  = "System, dnasr281@gmail.com"
    68  = "System, dnasr281@gmail.com"
    69  = "System, dnasr281@gmail.com"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    86  = "System, dnasr281@gmail.com"
    89  = "System, dnasr281@gmail.com"
    90  = "admin@admin.com, dnasr281@gmail.com"
    93  = "System, dnasr281@gmail.com"
    95  = "System, dnasr281@gmail.com"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    112 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    116 = "admin@admin.com, dnasr281@gmail.com"
    119 = "System, dnasr281@gmail.com"
    121 = "System, dnasr281@gmail.com"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    138 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    142 = "admin@admin.com, dnasr281@gmail.com"
    145 = "System, dnasr281@gmail.com"
    147 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
